$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Change the text in A18 from "RepaymentStrategy" to "repaymentstrategy" (lowercase)
$ws.Range("A18").Value = "repaymentstrategy"

# Update the selected cell/range to A18 (was B20)
$ws.Range("A18").Select()
